$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the existing hyperlink (was on F2 -> https://www.amazon.com/) before
#    we shift columns around, so we don't end up with stale relationship ids.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a brand new column before column A. Everything that used to live
#    in A:G now lives in B:H.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# ---------------------------------------------------------------------------
# 3. Give the new column A the same look as column B (which holds what used
#    to be column A) - bold/shaded header in row 1, plain wrapped text below.
# ---------------------------------------------------------------------------
$hdrSrc = $ws.Cells.Item(1, 2)
$hdrDst = $ws.Cells.Item(1, 1)
$hdrDst.Font.Bold = $hdrSrc.Font.Bold
$hdrDst.Font.Name = $hdrSrc.Font.Name
$hdrDst.Font.Size = $hdrSrc.Font.Size
$hdrDst.Interior.Color = $hdrSrc.Interior.Color
$hdrDst.Interior.Pattern = $hdrSrc.Interior.Pattern
$hdrDst.HorizontalAlignment = $hdrSrc.HorizontalAlignment
$hdrDst.VerticalAlignment = $hdrSrc.VerticalAlignment
$hdrDst.WrapText = $hdrSrc.WrapText

for ($r = 2; $r -le 11; $r++) {
    $bodySrc = $ws.Cells.Item($r, 2)
    $bodyDst = $ws.Cells.Item($r, 1)
    $bodyDst.Font.Bold = $bodySrc.Font.Bold
    $bodyDst.Font.Name = $bodySrc.Font.Name
    $bodyDst.Font.Size = $bodySrc.Font.Size
    $bodyDst.HorizontalAlignment = $bodySrc.HorizontalAlignment
    $bodyDst.VerticalAlignment = $bodySrc.VerticalAlignment
    $bodyDst.WrapText = $bodySrc.WrapText
}

# ---------------------------------------------------------------------------
# 4. Give the brand new column H (the old G, PATH TO IMAGES, shifted away)
#    the same look as column G - bold/shaded header in row 1, plain wrapped
#    text below for the data rows.
# ---------------------------------------------------------------------------
$hdrSrc2 = $ws.Cells.Item(1, 7)
$hdrDst2 = $ws.Cells.Item(1, 8)
$hdrDst2.Font.Bold = $hdrSrc2.Font.Bold
$hdrDst2.Font.Name = $hdrSrc2.Font.Name
$hdrDst2.Font.Size = $hdrSrc2.Font.Size
$hdrDst2.Interior.Color = $hdrSrc2.Interior.Color
$hdrDst2.Interior.Pattern = $hdrSrc2.Interior.Pattern
$hdrDst2.HorizontalAlignment = $hdrSrc2.HorizontalAlignment
$hdrDst2.VerticalAlignment = $hdrSrc2.VerticalAlignment
$hdrDst2.WrapText = $hdrSrc2.WrapText

# Row 3 never had a value as far right as column G/H in the source sheet, so
# it is intentionally skipped here - only row 2 and rows 4-11 get a new H.
foreach ($r in @(2, 4, 5, 6, 7, 8, 9, 10, 11)) {
    $bodySrc2 = $ws.Cells.Item($r, 2)
    $bodyDst2 = $ws.Cells.Item($r, 8)
    $bodyDst2.Font.Bold = $bodySrc2.Font.Bold
    $bodyDst2.Font.Name = $bodySrc2.Font.Name
    $bodyDst2.Font.Size = $bodySrc2.Font.Size
    $bodyDst2.HorizontalAlignment = $bodySrc2.HorizontalAlignment
    $bodyDst2.VerticalAlignment = $bodySrc2.VerticalAlignment
    $bodyDst2.WrapText = $bodySrc2.WrapText
}

# ---------------------------------------------------------------------------
# 5. New column A values - the "YOUR GMAIL" field.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "YOUR GMAIL"
$ws.Range("A2").Value = "user-example-mail@gmail.com"

# ---------------------------------------------------------------------------
# 6. Renamed headers (old EMAIL/PASSWORD/LINK columns, now C/D/G after the
#    insert).
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "PINTEREST EMAIL"
$ws.Range("D1").Value = "PINTEREST PASSWORD"
$ws.Range("G1").Value = "PINS LINK"

# ---------------------------------------------------------------------------
# 7. New link text + re-created hyperlink (now on G2). Add the hyperlink
#    first (Excel stamps its own built-in "Hyperlink" style when doing this),
#    then restore the workbook's own custom underlined-blue look so the cell
#    keeps matching the rest of the sheet's hand-rolled formatting.
# ---------------------------------------------------------------------------
$g2 = $ws.Range("G2")
$ws.Hyperlinks.Add($g2, "https://www.hugecakesexample.com/")
$g2.Value = "https://www.hugecakesexample.com/"
$g2.Font.Bold = $false
$g2.Font.Name = "Calibri"
$g2.Font.Size = 10
$g2.Font.Underline = 2
$g2.Font.Color = 16711680
$g2.HorizontalAlignment = 1
$g2.VerticalAlignment = -4160
$g2.WrapText = $true
